$d = $word.ActiveDocument
$lsq = [char]0x2018
$rsq = [char]0x2019

# ---------------------------------------------------------------------------
# Change 1 (paragraph "Select 'Install Software' from the Settings/System
# Configuration menu."): the three runs "Select ", "'Install Software'" and
# " from the Settings/" merge into a single run, while the following
# "System" and " Configuration menu." runs stay untouched/separate.
# ---------------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("Select " + $lsq + "Install Software" + $rsq + " from the Settings/")
$c1Start = $r.Start
$c1Text = $r.Text

# Force a real text mutation (append + trim) so the engine folds the three
# runs into one; then fix the text back to the exact target string.
$r.Text = $c1Text + " "
$r2 = $d.Range($c1Start, $r.End)
$r2.Text = $c1Text

# The merge above also absorbed "System"/" Configuration menu." into the same
# run; re-establish "System" as its own run (matching the unmodified source)
# by briefly wrapping it in a content control and removing the wrapper - this
# creates a run boundary without leaving any residual formatting behind.
$sysStart = $r2.End
$sysRange = $d.Range($sysStart, $sysStart + 6)
if ($sysRange.Text -eq "System") {
    $cc = $sysRange.ContentControls.Add()
    $cc.Delete()
}

# ---------------------------------------------------------------------------
# Change 2 (paragraph "The Install Software selection menu..."): fix the
# duplicated period before "click to install as required."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("not yet installed., click", $true, $false, $false, $false, $false, $true, 1, $false, "not yet installed, click", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3 (paragraph "The results of the last software installation task
# can be shown by selecting 'Install Software' from the Settings/System
# Status menu."): all six text runs merge into a single run.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$null = $r3.Find.Execute("The results of the last software installation task can be shown by selecting " + $lsq + "Install Software" + $rsq + " from the Settings/System Status menu.")
$c3Start = $r3.Start
$c3Text = $r3.Text

$r3.Text = $c3Text + " "
$r4 = $d.Range($c3Start, $r3.End)
$r4.Text = $c3Text
